# Tester för spara uppgift funkar
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsredovisning")
$ws.Activate()

# Insert a new row before the current "Summa" row (row 20), which pushes
# Summa down to row 21, then do it again so Summa ends at row 22.
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

# Fill in the two new data rows (20 and 21)
$ws.Range("A20").Value = 45314
$ws.Range("B20").Value = 2.5
$ws.Range("C20").Value = "Hämta uppgifter med sida + Test"

$ws.Range("A21").Value = 45315
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = "Hämta uppgifter med datum + Test"

# Update the Summa (totals) row formula to cover the new data range
$ws.Range("B22").Formula = "=SUBTOTAL(109,B2:B21)"

# Resize the table to include the two new rows (now A1:C22)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C22"))

# Widen column C
$ws.Columns.Item(3).ColumnWidth = 29.166666666666668

# Update selection to match the authored change
$ws.Range("D9").Select()

$excel.Calculate()
